# Update "想去人数" (F column) figures across the workbook's sheets
# to reflect newly generated output (gh-pages rebuild at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1821
$ws1.Range("F4").Value = 404
$ws1.Range("F5").Value = 1488
$ws1.Range("F6").Value = 854
$ws1.Range("F7").Value = 380
$ws1.Range("F8").Value = 738
$ws1.Range("F9").Value = 13180
$ws1.Range("F10").Value = 13046
$ws1.Range("F11").Value = 993
$ws1.Range("F16").Value = 637
$ws1.Range("F18").Value = 53
$ws1.Range("F19").Value = 28
$ws1.Range("F20").Value = 38
$ws1.Range("F22").Value = 192
$ws1.Range("F23").Value = 273

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 90
$ws2.Range("F7").Value = 105

# Sheet "本地生活" (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 188

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 188
$ws4.Range("F4").Value = 1821
$ws4.Range("F5").Value = 404
$ws4.Range("F6").Value = 1488
$ws4.Range("F7").Value = 854
$ws4.Range("F8").Value = 380
$ws4.Range("F9").Value = 90
$ws4.Range("F10").Value = 738
$ws4.Range("F11").Value = 13180
$ws4.Range("F12").Value = 13046
$ws4.Range("F13").Value = 993
$ws4.Range("F18").Value = 637
$ws4.Range("F22").Value = 53
$ws4.Range("F23").Value = 28
$ws4.Range("F24").Value = 38
$ws4.Range("F28").Value = 192
$ws4.Range("F29").Value = 273
$ws4.Range("F31").Value = 105
